$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '67.729.76'
$ws.Range("E2").Value = '  -0.59%  '

# Row 3
$ws.Range("D3").Value = '3.805.36'
$ws.Range("E3").Value = '  +1.81%  '

# Row 4
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").Value = '''595.83'
$ws.Range("D5").Style = $ws.Range("B5").Style
$ws.Range("E5").Value = '  +0.61%  '

# Row 6
$ws.Range("D6").Value = '''167.26'
$ws.Range("D6").Style = $ws.Range("B6").Style
$ws.Range("E6").Value = '  +0.21%  '

# Row 7
$ws.Range("D7").Value = '3.803.16'
$ws.Range("E7").Value = '  +1.81%  '

# Row 8
$ws.Range("E8").Value = '  +0.06%  '

# Row 9
$ws.Range("E9").Value = '  -0.03%  '

# Row 10
$ws.Range("D10").Value = '''0.160'
$ws.Range("D10").Style = $ws.Range("B10").Style
$ws.Range("E10").Value = '  +0.27%  '

# Row 11
$ws.Range("D11").Value = '''6.31'
$ws.Range("D11").Style = $ws.Range("B11").Style
$ws.Range("E11").Value = '  -1.55%  '

# Row 12
$ws.Range("D12").Value = '''0.449'
$ws.Range("D12").Style = $ws.Range("B12").Style
$ws.Range("E12").Value = '  +0.07%  '

# Row 13
$ws.Range("E13").Value = '  -1.56%  '

# Row 14
$ws.Range("D14").Value = '''36.13'
$ws.Range("D14").Style = $ws.Range("B14").Style

# Row 15
$ws.Range("D15").Value = '4.433.54'
$ws.Range("E15").Value = '  +1.57%  '

# Row 16
$ws.Range("D16").Value = '3.796.23'
$ws.Range("E16").Value = '  +2.48%  '

# Row 17
$ws.Range("D17").Value = '''18.47'
$ws.Range("D17").Style = $ws.Range("B17").Style
$ws.Range("E17").Value = '  +3.40%  '

# Row 18
$ws.Range("D18").Value = '67.725.42'
$ws.Range("E18").Value = '  -0.53%  '

# Row 19
$ws.Range("E19").Value = '  +0.45%  '

# Row 20
$ws.Range("E20").Value = '  -0.03%  '

# Row 21
$ws.Range("E21").Value = '  -6.01%  '

# Row 22
$ws.Range("D22").Value = '''460.15'
$ws.Range("D22").Style = $ws.Range("B22").Style
$ws.Range("E22").Value = '  -0.96%  '

# Row 23
$ws.Range("D23").Value = '''0.699'
$ws.Range("D23").Style = $ws.Range("B23").Style
$ws.Range("E23").Value = '  +0.48%  '

# Row 24
$ws.Range("D24").Value = '''0.0000157'
$ws.Range("D24").Style = $ws.Range("B24").Style
$ws.Range("E24").Value = '  +6.53%  '

# Row 25
$ws.Range("D25").Value = '''83.33'
$ws.Range("D25").Style = $ws.Range("B25").Style
$ws.Range("E25").Value = '  -0.49%  '

# Row 26
$ws.Range("D26").Value = '''12.05'
$ws.Range("D26").Style = $ws.Range("B26").Style
$ws.Range("E26").Value = '  +1.76%  '

# Row 27
$ws.Range("D27").Value = '''2.12'
$ws.Range("D27").Style = $ws.Range("B27").Style
$ws.Range("E27").Value = '  -2.37%  '

# Row 28
$ws.Range("D28").Value = '''10.02'
$ws.Range("D28").Style = $ws.Range("B28").Style
$ws.Range("E28").Value = '  +0.07%  '

# Row 29
$ws.Range("E29").Value = '  +0.13%  '

# Row 30
$ws.Range("E30").Value = '  +0.36%  '

# Row 31
$ws.Range("E31").Value = '  +4.11%  '

# Row 32
$ws.Range("D32").Value = '''7.25'
$ws.Range("D32").Style = $ws.Range("B32").Style
$ws.Range("E32").Value = '  -0.15%  '

# Row 33
$ws.Range("D33").Value = '''29.74'
$ws.Range("D33").Style = $ws.Range("B33").Style
$ws.Range("E33").Value = '  +0.03%  '

# Row 34
$ws.Range("D34").Value = '''1.00'
$ws.Range("D34").Style = $ws.Range("B34").Style
$ws.Range("E34").Value = '  -0.05%  '

# Row 35
$ws.Range("D35").Value = '''9.09'
$ws.Range("D35").Style = $ws.Range("B35").Style
$ws.Range("E35").Value = '  -0.17%  '

# Row 36
$ws.Range("D36").Value = '3.741.87'
$ws.Range("E36").Value = '  +1.36%  '

# Row 37
$ws.Range("E37").Value = '  -0.24%  '

# Row 38
$ws.Range("E38").Value = '  -1.19%  '

# Row 39
$ws.Range("D39").Value = '''0.138'
$ws.Range("D39").Style = $ws.Range("B39").Style
$ws.Range("E39").Value = '  +0.66%  '

# Row 40
$ws.Range("E40").Value = '  +0.23%  '

# Row 41
$ws.Range("D41").Value = '''5.78'
$ws.Range("D41").Style = $ws.Range("B41").Style
$ws.Range("E41").Value = '  +0.36%  '

# Row 42
$ws.Range("D42").Value = '''0.999'
$ws.Range("D42").Style = $ws.Range("B42").Style
$ws.Range("E42").Value = '  -0.07%  '

# Row 43
$ws.Range("E43").Value = '  +0.02%  '

# Row 44
$ws.Range("D44").Value = '''44.60'
$ws.Range("D44").Style = $ws.Range("B44").Style
$ws.Range("E44").Value = '  +1.16%  '

# Row 45
$ws.Range("D45").Value = '''48.05'
$ws.Range("D45").Style = $ws.Range("B45").Style
$ws.Range("E45").Value = '  +2.83%  '

# Row 46
$ws.Range("D46").Value = '''0.299'
$ws.Range("D46").Style = $ws.Range("B46").Style
$ws.Range("E46").Value = '  -0.11%  '

# Row 47
$ws.Range("D47").Value = '''149.00'
$ws.Range("D47").Style = $ws.Range("B47").Style
$ws.Range("E47").Value = '  +3.51%  '

# Row 48
$ws.Range("E48").Value = '  -1.34%  '

# Row 49
$ws.Range("D49").Value = '''394.67'
$ws.Range("D49").Style = $ws.Range("B49").Style
$ws.Range("E49").Value = '  +1.68%  '

# Row 50
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '''26.94'
$ws.Range("D50").Style = $ws.Range("B50").Style
$ws.Range("E50").Value = '  +7.36%  '

# Row 51
$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").Value = '''1.82'
$ws.Range("D51").Style = $ws.Range("B51").Style
$ws.Range("E51").Value = '  -3.99%  '
